$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("X5").Value = 28
$ws.Range("AD5").Value = 7.7
$ws.Range("AE5").Value = 12

# Row 6
$ws.Range("J6").Value = 1.17
$ws.Range("K6").Value = 5
$ws.Range("P6").Value = 1.62
$ws.Range("Q6").Value = 2.2
$ws.Range("R6").Value = 2.2
$ws.Range("S6").Value = 1.62
$ws.Range("T6").Value = 6.5
$ws.Range("Y6").Value = 41

# Row 9
$ws.Range("G9").Value = 2.32
$ws.Range("I9").Value = 3
$ws.Range("Q9").Value = 2.35

# Row 10
$ws.Range("J10").Value = 1.1
$ws.Range("K10").Value = 7
$ws.Range("L10").Value = 1.53
$ws.Range("M10").Value = 2.38

# Row 11
$ws.Range("H11").Value = 3.85
$ws.Range("S11").Value = 2.07
$ws.Range("T11").Value = 8.75
$ws.Range("U11").Value = 9
$ws.Range("W11").Value = 13.5
$ws.Range("X11").Value = 11.75
$ws.Range("Y11").Value = 20
$ws.Range("Z11").Value = 13.5
$ws.Range("AA11").Value = 7.7
$ws.Range("AD11").Value = 15
$ws.Range("AE11").Value = 29
$ws.Range("AI11").Value = 40

# Row 12
$ws.Range("H12").Value = 3.65
$ws.Range("U12").Value = 20
$ws.Range("Y12").Value = 27
$ws.Range("AA12").Value = 7.5
$ws.Range("AD12").Value = 10.25

# Row 13
$ws.Range("H13").Value = 3.3
$ws.Range("O13").Value = 1.9
$ws.Range("R13").Value = 1.57
$ws.Range("S13").Value = 2.1
$ws.Range("T13").Value = 9.25
$ws.Range("U13").Value = 12.5
$ws.Range("X13").Value = 16.5
$ws.Range("Y13").Value = 23
$ws.Range("AA13").Value = 6.5
$ws.Range("AD13").Value = 10.25
$ws.Range("AE13").Value = 16
$ws.Range("AF13").Value = 10.5
$ws.Range("AH13").Value = 24
$ws.Range("AI13").Value = 29
$ws.Range("AJ13").Value = 300

# Row 18
$ws.Range("J18").Value = 1.05
$ws.Range("K18").Value = 11
$ws.Range("N18").Value = 1.85
$ws.Range("O18").Value = 1.95

# Row 20
$ws.Range("K20").Value = 7.5
$ws.Range("P20").Value = 1.5
$ws.Range("Q20").Value = 2.5
$ws.Range("U20").Value = 8
$ws.Range("Z20").Value = 7.5

# Row 21
$ws.Range("I21").Value = 3.6
$ws.Range("K21").Value = 5.6
$ws.Range("N21").Value = 2.42
$ws.Range("P21").Value = 1.57
$ws.Range("Q21").Value = 2.25
$ws.Range("T21").Value = 5.9
$ws.Range("W21").Value = 20
$ws.Range("Y21").Value = 40
$ws.Range("Z21").Value = 5.6
$ws.Range("AD21").Value = 8
$ws.Range("AH21").Value = 40
$ws.Range("AI21").Value = 55

# Row 27
$ws.Range("G27").Value = 2.35
$ws.Range("I27").Value = 3
$ws.Range("R27").Value = 1.85
$ws.Range("S27").Value = 1.75
$ws.Range("U27").Value = 10.75
$ws.Range("V27").Value = 9.25
$ws.Range("W27").Value = 24
$ws.Range("X27").Value = 21
$ws.Range("AE27").Value = 14.5
$ws.Range("AF27").Value = 11
$ws.Range("AH27").Value = 30
$ws.Range("AI27").Value = 45

# Row 30
$ws.Range("G30").Value = 2.75
$ws.Range("I30").Value = 2.4
$ws.Range("N30").Value = 1.73
$ws.Range("O30").Value = 2.08
$ws.Range("W30").Value = 29
$ws.Range("X30").Value = 21
$ws.Range("AH30").Value = 17

# Row 31
$ws.Range("J31").Value = 1.03
$ws.Range("L31").Value = 1.25

# Row 32
$ws.Range("G32").Value = 2.75
$ws.Range("H32").Value = 3.4
$ws.Range("K32").Value = 7.9
$ws.Range("M32").Value = 3.6
$ws.Range("O32").Value = 1.98
$ws.Range("P32").Value = 1.37
$ws.Range("Q32").Value = 2.87
$ws.Range("S32").Value = 2.18
$ws.Range("Y32").Value = 27
$ws.Range("Z32").Value = 7.9
$ws.Range("AD32").Value = 9.25
$ws.Range("AI32").Value = 24

# Row 33
$ws.Range("I33").Value = 3.7
$ws.Range("R33").Value = 2.12
$ws.Range("S33").Value = 1.65
$ws.Range("T33").Value = 5.6
$ws.Range("U33").Value = 8.25
$ws.Range("W33").Value = 17.5
$ws.Range("X33").Value = 19.5
$ws.Range("Y33").Value = 40
$ws.Range("AA33").Value = 6.2
$ws.Range("AB33").Value = 19.5
$ws.Range("AF33").Value = 14
$ws.Range("AG33").Value = 60

# Row 34
$ws.Range("H34").Value = 4.35
$ws.Range("I34").Value = 6.4
$ws.Range("L34").Value = 1.23
$ws.Range("M34").Value = 3.75
$ws.Range("O34").Value = 2.05
$ws.Range("T34").Value = 7
$ws.Range("U34").Value = 6.7
$ws.Range("W34").Value = 9.25
$ws.Range("AD34").Value = 17.5
$ws.Range("AE34").Value = 40
$ws.Range("AH34").Value = 70
$ws.Range("AJ34").Value = 700

# Row 36
$ws.Range("K36").Value = 6.8
$ws.Range("P36").Value = 1.45
$ws.Range("Q36").Value = 2.55
$ws.Range("X36").Value = 18
$ws.Range("Z36").Value = 6.8
$ws.Range("AD36").Value = 9.25

# Row 38
$ws.Range("H38").Value = 3.4
$ws.Range("K38").Value = 7.4
$ws.Range("L38").Value = 1.3
$ws.Range("M38").Value = 3.2
$ws.Range("N38").Value = 1.9
$ws.Range("S38").Value = 1.95
$ws.Range("U38").Value = 11
$ws.Range("W38").Value = 22
$ws.Range("Z38").Value = 7.4
$ws.Range("AA38").Value = 6.6
$ws.Range("AB38").Value = 14.5
$ws.Range("AF38").Value = 10.25

# Row 39
$ws.Range("L39").Value = 1.33
$ws.Range("M39").Value = 3.05
$ws.Range("W39").Value = 37
$ws.Range("X39").Value = 27
$ws.Range("AD39").Value = 7.4

# Row 40
$ws.Range("K40").Value = 7
$ws.Range("L40").Value = 1.34
$ws.Range("M40").Value = 3
$ws.Range("N40").Value = 2
$ws.Range("O40").Value = 1.72
$ws.Range("Q40").Value = 2.6
$ws.Range("S40").Value = 1.88
$ws.Range("T40").Value = 7.3
$ws.Range("U40").Value = 10
$ws.Range("Z40").Value = 7
$ws.Range("AA40").Value = 6.5
$ws.Range("AB40").Value = 15.5

# Row 42
$ws.Range("G42").Value = 2.27
$ws.Range("H42").Value = 3.1
$ws.Range("I42").Value = 3.05
$ws.Range("L42").Value = 1.37
$ws.Range("M42").Value = 2.62
$ws.Range("N42").Value = 2.07
$ws.Range("O42").Value = 1.6
$ws.Range("P42").Value = 1.42
$ws.Range("Q42").Value = 2.47
$ws.Range("T42").Value = 6.7
$ws.Range("U42").Value = 10.25
$ws.Range("V42").Value = 9.25
$ws.Range("W42").Value = 23
$ws.Range("X42").Value = 21
$ws.Range("Z42").Value = 7.9
$ws.Range("AA42").Value = 6
$ws.Range("AD42").Value = 8.25
$ws.Range("AE42").Value = 15
$ws.Range("AF42").Value = 11
$ws.Range("AG42").Value = 40
$ws.Range("AH42").Value = 29
